# Apply weekly update to the Hortaliza / Pepino dulce sheet.
# The data rows (pairs/triples sharing the same Fecha) are rotated: each
# block's Fecha/Volumen/Precio.../Unidad values are replaced by the values
# that, before the edit, belonged to a different block - effectively
# refreshing the weekly price series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha (Excel serial date), Calidad, Volumen, PrecioMin, PrecioMax,
# PrecioProm, Unidad, PrecioKg, KgUnidades
$updates = @(
    @{ Row = 2;  D = 44756; I = "Primera"; J = 65; K = 14000; L = 14000; M = 14000; N = "`$/caja 15 kilos";     P = 933;  Q = 15 },
    @{ Row = 3;  D = 44756; I = "Segunda"; J = 68; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";     P = 800;  Q = 15 },
    @{ Row = 6;  D = 44238; I = "Primera"; J = 90; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos";  P = 722;  Q = 18 },
    @{ Row = 7;  D = 44238; I = "Segunda"; J = 80; K = 11000; L = 11000; M = 11000; N = "`$/bandeja 18 kilos";  P = 611;  Q = 18 },
    @{ Row = 8;  D = 44992; I = "Primera"; J = 56; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos";  P = 722;  Q = 18 },
    @{ Row = 9;  D = 44991; I = "Primera"; J = 75; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos";  P = 722;  Q = 18 },
    @{ Row = 10; D = 44991; I = "Segunda"; J = 56; K = 9000;  L = 9000;  M = 9000;  N = "`$/bandeja 18 kilos";  P = 500;  Q = 18 },
    @{ Row = 11; D = 44424; I = "Primera"; J = 75; K = 18000; L = 18000; M = 18000; N = "`$/caja 15 kilos";     P = 1200; Q = 15 },
    @{ Row = 12; D = 44424; I = "Segunda"; J = 50; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";     P = 800;  Q = 15 },
    @{ Row = 13; D = 44242; I = "Primera"; J = 60; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos";  P = 722;  Q = 18 },
    @{ Row = 14; D = 44242; I = "Segunda"; J = 50; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos";  P = 556;  Q = 18 },
    @{ Row = 15; D = 44235; I = "Primera"; J = 80; K = 14000; L = 14000; M = 14000; N = "`$/bandeja 18 kilos";  P = 778;  Q = 18 },
    @{ Row = 16; D = 44235; I = "Segunda"; J = 70; K = 12000; L = 12000; M = 12000; N = "`$/bandeja 18 kilos";  P = 667;  Q = 18 },
    @{ Row = 17; D = 44235; I = "Tercera"; J = 60; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos";  P = 556;  Q = 18 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D    # D - Fecha
    $ws.Cells.Item($r, 9).Value = $u.I    # I - Calidad
    $ws.Cells.Item($r, 10).Value = $u.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $u.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $u.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $u.N   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $u.P   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $u.Q   # Q - Kg o Unidades
}
